$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud'
$ws.Range("G3").Value = 'Dr. Manar Montaser, Dr. Alshimaa Atef, Dr. Gehan Adel, Administrator'
$ws.Range("G4").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Hanan Ragab, Dr. Heba Mahmoud Ali'
$ws.Range("G5").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab'
$ws.Range("G6").Value = 'Dr. Sara Nabil, Dr. Safa Hany'
$ws.Range("G9").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Marina Youhanna, Dr. Yasmeena Fattoh'
$ws.Range("G12").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G17").Value = 'Dr. Walaa Ghanima, Dr. Marian Samir, Dr. Enas Omran'
$ws.Range("G18").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G19").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Yasmin'
$ws.Range("G20").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Remon, Dr. Nardine, Dr. Monica'
$ws.Range("G21").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud'
$ws.Range("G22").Value = 'Dr. Manar Montaser, Dr. Alshimaa Atef, Dr. Gehan Adel, Administrator'
$ws.Range("G23").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Hanan Ragab, Dr. Heba Mahmoud Ali'
$ws.Range("G24").Value = 'Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Fatma Elhady'
$ws.Range("G25").Value = 'Dr. Yasmin Tarek, Dr. Nourhan Mohammad'
$ws.Range("G28").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Basma Hamed, Dr. Esraa Mostafa'
$ws.Range("G29").Value = 'Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Amira Ibrahim'
$ws.Range("G31").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G36").Value = 'Dr. Walaa Ghanima, Dr. Marian Samir, Dr. Enas Omran'
$ws.Range("G37").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G38").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Remon, Dr. Nardine, Dr. Monica'
$ws.Range("G39").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Remon, Dr. Nardine, Dr. Monica'
$ws.Range("G40").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud'
$ws.Range("G41").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud'
$ws.Range("G42").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G43").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Kerelos Zareef, Dr. Fatma Elhady'
$ws.Range("G44").Value = 'Dr. Sara Nabil, Dr. Safa Hany'
$ws.Range("G47").Value = 'Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Amira Ibrahim'
$ws.Range("G48").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Merna Said, Dr. Sarah Abdelmohsen, Dr. Yasmeena Fattoh, Dr. Fatma Shoukry'
$ws.Range("G49").Value = 'Dr. Mohammad Safwat, Dr. Mariam Toma Gerges'
$ws.Range("G50").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G56").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G57").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Remon, Dr. Nardine, Dr. Monica'
$ws.Range("G58").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Remon, Dr. Nardine, Dr. Monica'
$ws.Range("G59").Value = 'Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Nesma'
$ws.Range("G60").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud'
$ws.Range("G61").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Nahla Nagiub, Dr. Asmaa Reda'
$ws.Range("G63").Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Range("G66").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Marina Youhanna, Dr. Dina Adel, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim'
$ws.Range("G67").Value = 'Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Amira Ibrahim'
$ws.Range("G75").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G76").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Yasmin'
$ws.Range("G77").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Remon, Dr. Nardine, Dr. Monica'
$ws.Range("G78").Value = 'Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Nesma'
$ws.Range("G79").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud'
$ws.Range("G80").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Nahla Nagiub, Dr. Asmaa Reda'
$ws.Range("G81").Value = 'Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Fatma Elhady'
$ws.Range("G82").Value = 'Dr. Yasmin Tarek, Dr. Nourhan Mohammad'
$ws.Range("G83").Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Range("G85").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Marina Youhanna, Dr. Dina Adel, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim'
$ws.Range("G86").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Merna Said, Dr. Sarah Abdelmohsen, Dr. Yasmeena Fattoh, Dr. Fatma Shoukry'
$ws.Range("G88").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G94").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G95").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Yasmin'
$ws.Range("G96").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Remon, Dr. Nardine, Dr. Monica'
$ws.Range("G97").Value = 'Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Nesma'
$ws.Range("G98").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud'
$ws.Range("G99").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G100").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Kerelos Zareef, Dr. Fatma Elhady'
$ws.Range("G101").Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Range("G104").Value = 'Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Amira Ibrahim'
$ws.Range("G113").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G115").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Remon, Dr. Nardine, Dr. Monica'
